$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("son")

# ---------------------------------------------------------------------------
# Sheet1 ("Sheet1") data updates.
# Write order matters: it controls the order new strings are appended to the
# shared-string table, which must match the target file exactly.
# ---------------------------------------------------------------------------

# E2: account number, kept as text via a quote-prefix (creates the new
# "quote-prefix only" style used later for every quote-prefixed text cell).
$ws1.Cells.Item(2, 5).Value = "'045704070000966"

# G2 / H2: base64 password blob, duplicated into both columns.
$ws1.Cells.Item(2, 7).Value = "YWJjMTIz"
$ws1.Cells.Item(2, 8).Value = "YWJjMTIz"

# Header row (F1:I1) for the new password / pincode columns.
$ws1.Cells.Item(1, 7).Value = "passKHDN_nhap"
$ws1.Cells.Item(1, 8).Value = "passKHDN_duyet"
$ws1.Cells.Item(1, 6).Value = "passKHCN"
$ws1.Cells.Item(1, 9).Value = "pincode"

# I2: new pincode column value.
$ws1.Cells.Item(2, 9).Value = "MTEyMjMz"

# C2 / D2: fill in the autotest user columns (values already exist elsewhere
# in the shared-string table, so this does not introduce new entries).
$ws1.Cells.Item(2, 3).Value = "autotestnhap"
$ws1.Cells.Item(2, 4).Value = "autotestduyet"

# Columns G:H widen to match column F.
$ws1.Columns.Item(7).ColumnWidth = $ws1.Columns.Item(6).ColumnWidth
$ws1.Columns.Item(8).ColumnWidth = $ws1.Columns.Item(6).ColumnWidth

$ws1.Activate()
$ws1.Range("E12").Select()

# ---------------------------------------------------------------------------
# Sheet2 ("son") data updates - reuses strings already introduced above, so
# the shared-string table append order is unaffected by what follows.
# ---------------------------------------------------------------------------

$ws2.Cells.Item(1, 6).Value = "passKHCN"
$ws2.Cells.Item(1, 7).Value = "passKHDN_nhap"
$ws2.Cells.Item(1, 8).Value = "passKHDN_duyet"

$ws2.Cells.Item(2, 1).Value = "minhson0907"

# B2 needs the same style as Sheet1!B2: quote-prefixed "0" number format.
$ws2.Cells.Item(2, 2).Value = "'002704070016025"
$ws2.Cells.Item(2, 2).NumberFormat = "0"

$ws2.Cells.Item(2, 5).Value = "'045704070000966"

$ws2.Cells.Item(2, 6).Value = "YXBwbGVpcGhvbmU2Uw=="
$ws2.Cells.Item(2, 7).Value = "YWJjMTIz"
$ws2.Cells.Item(2, 8).Value = "YWJjMTIz"

$ws2.Activate()
$ws2.Range("E33").Select()

$ws1.Activate()
